{"js": "// Remove the \"Maximal Square\" entry (and the blank separator paragraph that\n// follows it) from the DP section, and remove the \"Word Ladder II\" entry\n// from the end of the second Arrays section.\n\n// 1) \"Maximal Square\" + its trailing blank paragraph.\nconst msResults = context.document.body.search(\"Maximal Square\", { matchCase: true });\nmsResults.load(\"items\");\nawait context.sync();\n\nif (msResults.items.length > 0) {\n  const msPara = msResults.items[0].paragraphs.getFirst();\n  const afterPara = msPara.getNext();\n  afterPara.delete();\n  msPara.delete();\n}\n\n// 2) \"Word Ladder II\".\nconst wlResults = context.document.body.search(\"Word Ladder II\", { matchCase: true });\nwlResults.load(\"items\");\nawait context.sync();\n\nif (wlResults.items.length > 0) {\n  const wlPara = wlResults.items[0].paragraphs.getFirst();\n  wlPara.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Maximal Square\" entry (and the blank separator paragraph that\n# follows it) from the DP section, and remove the \"Word Ladder II\" entry\n# from the end of the second Arrays section.\n\n$d = $word.ActiveDocument\n\n# 1) \"Maximal Square\" + its trailing blank paragraph.\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute(\"Maximal Square\")\nif ($found) {\n    [void]$rng.Expand(4)              # wdParagraph - include the paragraph mark\n    $nextPara = $d.Range($rng.End, $rng.End)\n    [void]$nextPara.Expand(4)         # the blank paragraph right after it\n    $delRange = $d.Range($rng.Start, $nextPara.End)\n    [void]$delRange.Delete()\n}\n\n# 2) \"Word Ladder II\".\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$found2 = $rng2.Find.Execute(\"Word Ladder II\")\nif ($found2) {\n    [void]$rng2.Expand(4)\n    [void]$rng2.Delete()\n}\n"}
